$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Columns: 1=Watershed, 2=LF, 3=Rank, 4=Total Risk, 5=Current Risk, 6=Future Risk

# Row 10 (LF65: Mortality or fitness reduction due to deleterious substances)
# Total Risk 6 -> 8, Future Risk M -> H
$t.Cell(10, 4).Range.Text = "8"
$t.Cell(10, 6).Range.Text = "H"

# Row 11 (LF40: Mortality or fitness reduction due to frequent and higher peak flows causing flushing)
# Total Risk 4 -> 6, Future Risk L -> M
$t.Cell(11, 4).Range.Text = "6"
$t.Cell(11, 6).Range.Text = "M"

# Row 12 (was LF6) becomes LF9 row: LF text, Total Risk, Current Risk change; Rank/Future stay same
$t.Cell(12, 2).Range.Text = "LF9: Mortality or fitness reduction due to fishing"
$t.Cell(12, 4).Range.Text = "4"
$t.Cell(12, 5).Range.Text = "L"

# Row 13 (was LF9) becomes LF6 row: LF text, Rank, Total Risk, Current Risk, Future Risk change
$t.Cell(13, 2).Range.Text = "LF6: Limited or delayed access due to physical migration barriers and/or lack of safe migration routes (including lack of cover and complexity)"
$t.Cell(13, 3).Range.Text = "12"
$t.Cell(13, 4).Range.Text = "3"
$t.Cell(13, 5).Range.Text = "VL"
$t.Cell(13, 6).Range.Text = "M"

# Row 14 (LF56: Mortality or fitness reduction due to reduction in quality channel habitat)
# Rank 11 -> 13
$t.Cell(14, 3).Range.Text = "13"

# Row 15 (LF57: Mortality or fitness reduction due to reduction in quantity channel habitat)
# Rank 11 -> 13
$t.Cell(15, 3).Range.Text = "13"
